$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Headers: BTec_Logo-Orange pictures rename image1.jpg -> image2.jpg
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
            $shp = $hdr.Range.InlineShapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}

# Footers: PearsonLogo pictures rename image2.png -> image1.png
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
            $shp = $ftr.Range.InlineShapes.Item($j)
            if ($shp.AlternativeText -like "*PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}

Write-Output "done"
